$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.033.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "'3.570.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'602.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").Value = "'135.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D7").Value = "'3.566.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.58%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "'4.175.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "'3.569.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "'65.114.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "'10.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.91%  "
$ws.Range("D20").Value = "'14.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.44%  "
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").Value = "'388.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'0.580"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.63%  "
$ws.Range("D24").Value = "'3.715.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").Value = "'74.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'0.0000116"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.38%  "
$ws.Range("D28").Value = "'7.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.68%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +4.31%  "
$ws.Range("D31").Value = "'8.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("E32").Value = "  +21.68%  "
$ws.Range("D33").Value = "'3.578.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").Value = "'24.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.52%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").Value = "'6.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "'169.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "'1.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.66%  "
$ws.Range("D40").Value = "'5.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.30%  "
$ws.Range("D41").Value = "'0.0808"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.37%  "
$ws.Range("D42").Value = "'27.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.63%  "
$ws.Range("D43").Value = "'0.828"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'42.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'4.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.12%  "
$ws.Range("D47").Value = "'1.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.51%  "
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("D49").Value = "'2.499.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.57%  "
$ws.Range("E50").Value = "  +4.35%  "
$ws.Range("D51").Value = "'2.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.70%  "
